$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.068598123028721
$ws.Cells.Item(2, 4).Value = 1.069214107346037
$ws.Cells.Item(2, 5).Value = 1.072716216237676
$ws.Cells.Item(2, 6).Value = 1.081869338050028
$ws.Cells.Item(2, 9).Value = 1.053098109137426
$ws.Cells.Item(2, 10).Value = 1.073537116518704
$ws.Cells.Item(2, 11).Value = 1.07191721343545
$ws.Cells.Item(2, 12).Value = 1.075409992178965
$ws.Cells.Item(2, 13).Value = 1.084539036760271
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.069752685081887
$ws.Cells.Item(3, 4).Value = 1.070115331736703
$ws.Cells.Item(3, 5).Value = 1.073732100620998
$ws.Cells.Item(3, 6).Value = 1.082935817633668
$ws.Cells.Item(3, 9).Value = 1.053430970584871
$ws.Cells.Item(3, 10).Value = 1.07434753643135
$ws.Cells.Item(3, 11).Value = 1.072634353799943
$ws.Cells.Item(3, 12).Value = 1.076242183385786
$ws.Cells.Item(3, 13).Value = 1.085423441781775
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.070499949551975
$ws.Cells.Item(4, 4).Value = 1.070698609005514
$ws.Cells.Item(4, 5).Value = 1.074389896539056
$ws.Cells.Item(4, 6).Value = 1.083626430059324
$ws.Cells.Item(4, 9).Value = 1.053645243351831
$ws.Cells.Item(4, 10).Value = 1.074871529106606
$ws.Cells.Item(4, 11).Value = 1.073097879451308
$ws.Cells.Item(4, 12).Value = 1.076780502638637
$ws.Cells.Item(4, 13).Value = 1.0859956330061
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.070814144591256
$ws.Cells.Item(5, 4).Value = 1.070943848465059
$ws.Cells.Item(5, 5).Value = 1.074666541792544
$ws.Cells.Item(5, 6).Value = 1.083916890066406
$ws.Cells.Item(5, 9).Value = 1.053735057681415
$ws.Cells.Item(5, 10).Value = 1.075091719204988
$ws.Cells.Item(5, 11).Value = 1.073292623020768
$ws.Cells.Item(5, 12).Value = 1.077006772778794
$ws.Cells.Item(5, 13).Value = 1.086236163523281
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.070866901981097
$ws.Cells.Item(6, 4).Value = 1.07098502698538
$ws.Cells.Item(6, 5).Value = 1.074712998070246
$ws.Cells.Item(6, 6).Value = 1.083965666998633
$ws.Cells.Item(6, 9).Value = 1.053750122307586
$ws.Cells.Item(6, 10).Value = 1.075128684440243
$ws.Cells.Item(6, 11).Value = 1.073325314132525
$ws.Cells.Item(6, 12).Value = 1.077044762228377
$ws.Cells.Item(6, 13).Value = 1.086276548560274
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.070504147663639
$ws.Cells.Item(7, 4).Value = 1.070701885790724
$ws.Cells.Item(7, 5).Value = 1.074393592662782
$ws.Cells.Item(7, 6).Value = 1.083630310702992
$ws.Cells.Item(7, 9).Value = 1.053646444500265
$ws.Cells.Item(7, 10).Value = 1.074874471676215
$ws.Cells.Item(7, 11).Value = 1.073100482106217
$ws.Cells.Item(7, 12).Value = 1.076783526225622
$ws.Cells.Item(7, 13).Value = 1.085998847060025
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.068988274617784
$ws.Cells.Item(8, 4).Value = 1.069518654389779
$ws.Cells.Item(8, 5).Value = 1.073059445739118
$ws.Cells.Item(8, 6).Value = 1.082229650287651
$ws.Cells.Item(8, 9).Value = 1.053210831228346
$ws.Cells.Item(8, 10).Value = 1.073811085270119
$ws.Cells.Item(8, 11).Value = 1.072159680409148
$ws.Cells.Item(8, 12).Value = 1.075691269080399
$ws.Cells.Item(8, 13).Value = 1.084837941715139
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.066318506330669
$ws.Cells.Item(9, 4).Value = 1.067434613358311
$ws.Cells.Item(9, 5).Value = 1.070711963607995
$ws.Cells.Item(9, 6).Value = 1.079765559901797
$ws.Cells.Item(9, 9).Value = 1.052434717154804
$ws.Cells.Item(9, 10).Value = 1.071934169495963
$ws.Cells.Item(9, 11).Value = 1.070497945570388
$ws.Cells.Item(9, 12).Value = 1.073765314589346
$ws.Cells.Item(9, 13).Value = 1.08279168104341
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.064539556177625
$ws.Cells.Item(10, 4).Value = 1.066045908775603
$ws.Cells.Item(10, 5).Value = 1.069149299354002
$ws.Cells.Item(10, 6).Value = 1.078125557279785
$ws.Cells.Item(10, 9).Value = 1.051911586127511
$ws.Cells.Item(10, 10).Value = 1.070680803198111
$ws.Cells.Item(10, 11).Value = 1.069387483238042
$ws.Cells.Item(10, 12).Value = 1.072480495256541
$ws.Cells.Item(10, 13).Value = 1.081427103033606
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.063769453377603
$ws.Cells.Item(11, 4).Value = 1.065444740842335
$ws.Cells.Item(11, 5).Value = 1.068473198884212
$ws.Cells.Item(11, 6).Value = 1.077416063191965
$ws.Cells.Item(11, 9).Value = 1.05168370674739
$ws.Cells.Item(11, 10).Value = 1.070137583336886
$ws.Cells.Item(11, 11).Value = 1.068906013591176
$ws.Cells.Item(11, 12).Value = 1.071923951071887
$ws.Cells.Item(11, 13).Value = 1.080836128338913
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.06348343113237
$ws.Cells.Item(12, 4).Value = 1.065221462857187
$ws.Cells.Item(12, 5).Value = 1.068222146543281
$ws.Cells.Item(12, 6).Value = 1.077152621134449
$ws.Cells.Item(12, 9).Value = 1.051598857713459
$ws.Cells.Item(12, 10).Value = 1.069935731226748
$ws.Cells.Item(12, 11).Value = 1.06872707915533
$ws.Cells.Item(12, 12).Value = 1.07171719414725
$ws.Cells.Item(12, 13).Value = 1.080616598258584
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.06354478261336
$ws.Cells.Item(13, 4).Value = 1.065269355722893
$ws.Cells.Item(13, 5).Value = 1.068275994438387
$ws.Cells.Item(13, 6).Value = 1.077209126025423
$ws.Cells.Item(13, 9).Value = 1.051617067374778
$ws.Cells.Item(13, 10).Value = 1.069979032637344
$ws.Cells.Item(13, 11).Value = 1.068765465509159
$ws.Cells.Item(13, 12).Value = 1.071761545644865
$ws.Cells.Item(13, 13).Value = 1.080663688923118
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.063745810119465
$ws.Cells.Item(14, 4).Value = 1.065426284154547
$ws.Cells.Item(14, 5).Value = 1.068452445158282
$ws.Cells.Item(14, 6).Value = 1.077394285039472
$ws.Cells.Item(14, 9).Value = 1.051676697271981
$ws.Cells.Item(14, 10).Value = 1.070120899732244
$ws.Cells.Item(14, 11).Value = 1.068891224757519
$ws.Cells.Item(14, 12).Value = 1.071906861126735
$ws.Cells.Item(14, 13).Value = 1.080817982236921
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.06386967351149
$ws.Cells.Item(15, 4).Value = 1.065522975928065
$ws.Cells.Item(15, 5).Value = 1.068561173061768
$ws.Cells.Item(15, 6).Value = 1.077508380286039
$ws.Cells.Item(15, 9).Value = 1.051713410119415
$ws.Cells.Item(15, 10).Value = 1.070208298641626
$ws.Cells.Item(15, 11).Value = 1.068968696563868
$ws.Cells.Item(15, 12).Value = 1.071996390592004
$ws.Cells.Item(15, 13).Value = 1.080913045343965
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.064590669375618
$ws.Cells.Item(16, 4).Value = 1.066085809475757
$ws.Cells.Item(16, 5).Value = 1.069194181350922
$ws.Cells.Item(16, 6).Value = 1.078172657488172
$ws.Cells.Item(16, 9).Value = 1.051926681062425
$ws.Cells.Item(16, 10).Value = 1.070716844263061
$ws.Cells.Item(16, 11).Value = 1.06941942345737
$ws.Cells.Item(16, 12).Value = 1.072517426838241
$ws.Cells.Item(16, 13).Value = 1.081466321879405
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.06504298237719
$ws.Cells.Item(17, 4).Value = 1.066438900312277
$ws.Cells.Item(17, 5).Value = 1.069591396105844
$ws.Cells.Item(17, 6).Value = 1.078589511664597
$ws.Cells.Item(17, 9).Value = 1.0520600960342
$ws.Cells.Item(17, 10).Value = 1.071035706318366
$ws.Cells.Item(17, 11).Value = 1.069701983148324
$ws.Cells.Item(17, 12).Value = 1.072844202999432
$ws.Cells.Item(17, 13).Value = 1.081813349784479
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.065306827716249
$ws.Cells.Item(18, 4).Value = 1.066644866728696
$ws.Cells.Item(18, 5).Value = 1.069823137199999
$ws.Cells.Item(18, 6).Value = 1.078832717211376
$ws.Cells.Item(18, 9).Value = 1.052137783499864
$ws.Cells.Item(18, 10).Value = 1.071221644556739
$ws.Cells.Item(18, 11).Value = 1.069866734457671
$ws.Cells.Item(18, 12).Value = 1.07303478586275
$ws.Cells.Item(18, 13).Value = 1.082015755327844
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.065396795380063
$ws.Cells.Item(19, 4).Value = 1.066715098415335
$ws.Cells.Item(19, 5).Value = 1.069902163825715
$ws.Cells.Item(19, 6).Value = 1.078915654465154
$ws.Cells.Item(19, 9).Value = 1.052164250663171
$ws.Cells.Item(19, 10).Value = 1.071285036457551
$ws.Cells.Item(19, 11).Value = 1.069922900059846
$ws.Cells.Item(19, 12).Value = 1.073099766284904
$ws.Cells.Item(19, 13).Value = 1.082084768704693
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.064994451533234
$ws.Cells.Item(20, 4).Value = 1.066401015506899
$ws.Cells.Item(20, 5).Value = 1.069548773298029
$ws.Cells.Item(20, 6).Value = 1.078544780790559
$ws.Cells.Item(20, 9).Value = 1.052045795445004
$ws.Cells.Item(20, 10).Value = 1.07100150046723
$ws.Cells.Item(20, 11).Value = 1.069671673476179
$ws.Cells.Item(20, 12).Value = 1.072809145095038
$ws.Cells.Item(20, 13).Value = 1.081776118006501
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.063686611758109
$ws.Cells.Item(21, 4).Value = 1.065380072020457
$ws.Cells.Item(21, 5).Value = 1.068400482571409
$ws.Cells.Item(21, 6).Value = 1.077339757686321
$ws.Cells.Item(21, 9).Value = 1.051659143399351
$ws.Cells.Item(21, 10).Value = 1.070079125511349
$ws.Cells.Item(21, 11).Value = 1.068854194422154
$ws.Cells.Item(21, 12).Value = 1.071864070214434
$ws.Cells.Item(21, 13).Value = 1.080772547138706
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.062864483795813
$ws.Cells.Item(22, 4).Value = 1.064738294794951
$ws.Cells.Item(22, 5).Value = 1.067678977931867
$ws.Cells.Item(22, 6).Value = 1.076582665063284
$ws.Cells.Item(22, 9).Value = 1.051414856756399
$ws.Cells.Item(22, 10).Value = 1.069498751497279
$ws.Cells.Item(22, 11).Value = 1.068339662163456
$ws.Cells.Item(22, 12).Value = 1.071269680980906
$ws.Cells.Item(22, 13).Value = 1.080141470975537
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.06330029427609
$ws.Cells.Item(23, 4).Value = 1.065078500688652
$ws.Cells.Item(23, 5).Value = 1.068061416574682
$ws.Cells.Item(23, 6).Value = 1.076983961812737
$ws.Cells.Item(23, 9).Value = 1.051544469920947
$ws.Cells.Item(23, 10).Value = 1.069806460591793
$ws.Cells.Item(23, 11).Value = 1.068612477718944
$ws.Cells.Item(23, 12).Value = 1.071584795401161
$ws.Cells.Item(23, 13).Value = 1.080476025122901
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.065016380474159
$ws.Cells.Item(24, 4).Value = 1.066418133974766
$ws.Cells.Item(24, 5).Value = 1.069568032547983
$ws.Cells.Item(24, 6).Value = 1.07856499255674
$ws.Cells.Item(24, 9).Value = 1.052052257671086
$ws.Cells.Item(24, 10).Value = 1.071016956768974
$ws.Cells.Item(24, 11).Value = 1.069685369300617
$ws.Cells.Item(24, 12).Value = 1.072824986315054
$ws.Cells.Item(24, 13).Value = 1.081792941474958
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.067008544571608
$ws.Cells.Item(25, 4).Value = 1.067973272865025
$ws.Cells.Item(25, 5).Value = 1.071318434472575
$ws.Cells.Item(25, 6).Value = 1.080402106142702
$ws.Cells.Item(25, 9).Value = 1.052636369625017
$ws.Cells.Item(25, 10).Value = 1.072419764741992
$ws.Cells.Item(25, 11).Value = 1.07092800886474
$ws.Cells.Item(25, 12).Value = 1.074263369463478
$ws.Cells.Item(25, 13).Value = 1.08332075952549
